# Trade #83 closed at 2026-02-18 00:28:51 - unknown UNKNOWN +0.000%
#
# Applies:
#  - Summary sheet: refreshed aggregate metrics (capital, P&L, trade counts, win rate)
#  - Strategy Status sheet: refreshed "momentum" strategy row
#  - All Trades sheet: trade #111 (momentum) closed out + two new open trades appended
#  - momentum sheet: trade #111 closed out (mirrors "All Trades")
#  - HighProbConvergence sheet: new open trade #140 appended
#  - MarketMaking sheet: new open trade #141 appended

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.41   # Current Capital
$summary.Range("B4").Value = 0.52      # Total P&L $
$summary.Range("B5").Value = 0.09      # Total P&L %
$summary.Range("B6").Value = 111       # Total Trades
$summary.Range("B8").Value = 41        # Losing Trades
$summary.Range("B9").Value = 47.75     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - "momentum" row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.23999999999999
$status.Range("D11").Value = 29
$status.Range("E11").Value = -0.75
$status.Range("F11").Value = -0.76
$status.Range("G11").Value = 27.59

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #111 (row 112) closes out
$allTrades.Range("G112").Value = 0.949591
$allTrades.Range("H112").Value = "CLOSED"
$allTrades.Range("I112").Value = -2.104
$allTrades.Range("J112").Value = -0.02
$allTrades.Range("K112").Value = 99.23999999999999
$allTrades.Range("L112").Value = "early_exit"
$allTrades.Range("M112").Value = 0.14

# New trade row 141 -> Trade #140 (HighProbConvergence, still OPEN)
$allTrades.Range("A141").Value = 140
$allTrades.Range("B141").NumberFormat = "@"
$allTrades.Range("B141").Value = "2026-02-18"
$allTrades.Range("C141").NumberFormat = "@"
$allTrades.Range("C141").Value = "00:28:44"
$allTrades.Range("D141").Value = "HighProbConvergence"
$allTrades.Range("E141").Value = "DOWN"
$allTrades.Range("F141").Value = 0.97
$allTrades.Range("H141").Value = "OPEN"
$allTrades.Range("I141").Value = 0
$allTrades.Range("J141").Value = 0
$allTrades.Range("K141").Value = 100.4130057263667
$allTrades.Range("M141").Value = 0
$allTrades.Range("N141").Value = 0
$allTrades.Range("O141").Value = 0
$allTrades.Range("P141").Value = 0.95
$allTrades.Range("Q141").Value = "Mean reversion DOWN: price 2.49% above mean (z=2.38)"

# New trade row 142 -> Trade #141 (MarketMaking, still OPEN)
$allTrades.Range("A142").Value = 141
$allTrades.Range("B142").NumberFormat = "@"
$allTrades.Range("B142").Value = "2026-02-18"
$allTrades.Range("C142").NumberFormat = "@"
$allTrades.Range("C142").Value = "00:28:45"
$allTrades.Range("D142").Value = "MarketMaking"
$allTrades.Range("E142").Value = "DOWN"
$allTrades.Range("F142").Value = 0.97
$allTrades.Range("H142").Value = "OPEN"
$allTrades.Range("I142").Value = 0
$allTrades.Range("J142").Value = 0
$allTrades.Range("K142").Value = 99.47967800952271
$allTrades.Range("M142").Value = 0
$allTrades.Range("N142").Value = 0
$allTrades.Range("O142").Value = 0
$allTrades.Range("P142").Value = 0.6
$allTrades.Range("Q142").Value = "Normal spread capture: 190 bps"

# ---------------------------------------------------------------------------
# momentum sheet - trade #111 (row 30) closes out, mirrors "All Trades"
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("G30").Value = 0.949591
$momentum.Range("H30").Value = "CLOSED"
$momentum.Range("I30").Value = -2.104
$momentum.Range("J30").Value = -0.02
$momentum.Range("K30").Value = 99.23999999999999
$momentum.Range("P30").Value = "early_exit"
$momentum.Range("Q30").Value = 0.14

# ---------------------------------------------------------------------------
# HighProbConvergence sheet - new trade row 16 -> Trade #140 (still OPEN)
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Range("A16").Value = 140
$hpc.Range("B16").NumberFormat = "@"
$hpc.Range("B16").Value = "2026-02-18"
$hpc.Range("C16").NumberFormat = "@"
$hpc.Range("C16").Value = "00:28:44"
$hpc.Range("D16").Value = "HighProbConvergence"
$hpc.Range("E16").Value = "DOWN"
$hpc.Range("F16").Value = 0.97
$hpc.Range("H16").Value = "OPEN"
$hpc.Range("I16").Value = 0
$hpc.Range("J16").Value = 0
$hpc.Range("K16").Value = 100.4130057263667
$hpc.Range("L16").Value = 0
$hpc.Range("M16").Value = 0
$hpc.Range("N16").Value = 0.95
$hpc.Range("O16").Value = "Mean reversion DOWN: price 2.49% above mean (z=2.38)"
$hpc.Range("Q16").Value = 0

# ---------------------------------------------------------------------------
# MarketMaking sheet - new trade row 58 -> Trade #141 (still OPEN)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("A58").Value = 141
$mm.Range("B58").NumberFormat = "@"
$mm.Range("B58").Value = "2026-02-18"
$mm.Range("C58").NumberFormat = "@"
$mm.Range("C58").Value = "00:28:45"
$mm.Range("D58").Value = "MarketMaking"
$mm.Range("E58").Value = "DOWN"
$mm.Range("F58").Value = 0.97
$mm.Range("H58").Value = "OPEN"
$mm.Range("I58").Value = 0
$mm.Range("J58").Value = 0
$mm.Range("K58").Value = 99.47967800952271
$mm.Range("L58").Value = 0
$mm.Range("M58").Value = 0
$mm.Range("N58").Value = 0.6
$mm.Range("O58").Value = "Normal spread capture: 190 bps"
$mm.Range("Q58").Value = 0
